$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, styled like existing headers (copy format from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-11
$values = @(
    @(8, 9),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(5, 5),
    @(8, 8)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $rowNum = $r + 2
    $ws.Cells.Item($rowNum, 9).Value = $values[$r][0]
    $ws.Cells.Item($rowNum, 10).Value = $values[$r][1]
}
